# GANTT A JOUR 27/10
# Advance the scroll-increment (cell E3, linked to the "Barre de défilement"
# scrollbar form control) from 25 to 30 days, and mark the
# "Réinitialisation par envoi de mail..." milestone (row 18) as 100% done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")

# Scroll-increment cell driving the whole Gantt timeline (named range
# Incrément_Défilement). Bumping it shifts the visible date window.
$ws.Range("E3").Value = 30

# Keep the scrollbar form control (linked to $E$3) in sync with the new value.
$scrollBar = $ws.Shapes.Item(1)
$scrollBar.ControlFormat.Value = 30

# Milestone progress update: "Avancement" 80% -> 100%.
$ws.Range("D18").Value = 1

# Leave the selection where the author left off.
$null = $ws.Range("E23").Select()
